$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Use")

$values = @{
    106 = "RoW"
    107 = "RoW"
    108 = "RoW"
    109 = "RoW"
    110 = "RoW"
    111 = "RoW"
    112 = "RoW"
    113 = "RoW"
    114 = "RoW"
    115 = "RoW"
    116 = "RoW"
    117 = "RoW"
    118 = "RoW"
    119 = "EU27"
    120 = "EU27"
    121 = "EU27"
    122 = "EU27"
    123 = "US"
    124 = "US"
    125 = "US"
    126 = "US"
    127 = "EU27"
    128 = "EU27"
    129 = "US"
    130 = "US"
    131 = "RoW"
    132 = "EU27"
    133 = "EU27"
    134 = "EU27"
    135 = "EU27"
    136 = "US"
    137 = "US"
    138 = "US"
    139 = "US"
    140 = "EU27"
    141 = "EU27"
    142 = "US"
    143 = "US"
    144 = "RoW"
    145 = "EU27"
    146 = "EU27"
    147 = "EU27"
    148 = "EU27"
    149 = "US"
    150 = "US"
    151 = "US"
    152 = "US"
    153 = "EU27"
    154 = "EU27"
    155 = "US"
    156 = "US"
    157 = "RoW"
    158 = "EU27"
    159 = "EU27"
    160 = "EU27"
    161 = "EU27"
    162 = "US"
    163 = "US"
    164 = "US"
    165 = "US"
    166 = "EU27"
    167 = "EU27"
    168 = "US"
    169 = "US"
    170 = "RoW"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value2 = $values[$row]
}

$ws.Range("C35").Select()

$wb.Save()